$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> new value for column F ("dSF")
$updates = @{
    2  = 6
    6  = -6
    7  = 2
    10 = -2
    13 = 2
    14 = -7
    15 = 0
    16 = -1
    18 = 1
    19 = 1
    21 = 0
    22 = 4
    25 = 5
    29 = -2
    30 = 2
    32 = 2
    33 = -1
    36 = -1
    37 = -1
    41 = -2
    43 = 0
    44 = -1
    46 = 4
    50 = -2
    51 = -5
    53 = -5
    57 = 0
    60 = -5
    61 = -5
    63 = -4
    64 = 1
    65 = -5
    69 = 5
    71 = -3
    77 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
